# library_3275.xlsx — add 3 new rows of library-prep records (rows 38-40),
# highlighting the (re)used index sequences in purple, matching the
# "corrected 37C.CO2 to DMEM.37C.C02" re-check of the index assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Purple highlight color used to flag re-used/duplicate index sequences
# (fgColor FFB767FF as solid fill).
$purple = 16738231   # 0xFFB767FF -> BGR 0xFF67B7 == RGB(183,103,255)

# --- Row 38 --------------------------------------------------------------
$ws.Range("A38").Value = "08.14.18"
$ws.Range("B38").Value = "H.BROWN"
$ws.Range("C38").Value = 38
$ws.Range("D38").Value = "08.14.18"
$ws.Range("E38").Value = "H.BROWN"
$ws.Range("F38").Value = 38
$ws.Range("G38").Value = "ATCGAGC"
$ws.Range("H38").Value = "E7420L"

$ws.Range("G38").Interior.Color = $purple
$ws.Range("H38").WrapText = $true

$ws.Rows.Item(38).RowHeight = 16

# --- Row 39 --------------------------------------------------------------
$ws.Range("A39").Value = "10.18.18"
$ws.Range("B39").Value = "H.BROWN"
$ws.Range("C39").Value = 39
$ws.Range("D39").Value = "10.18.18"
$ws.Range("E39").Value = "H.BROWN"
$ws.Range("F39").Value = 39
$ws.Range("G39").Value = "ATCGAGC"
$ws.Range("H39").Value = "E7420L"

$ws.Range("A39").Interior.Color = $purple
$ws.Range("D39").Interior.Color = $purple
$ws.Range("G39").Interior.Color = $purple
$ws.Range("H39").WrapText = $true

$ws.Rows.Item(39).RowHeight = 16

# --- Row 40 --------------------------------------------------------------
$ws.Range("A40").Value = "11.02.18"
$ws.Range("B40").Value = "H.BROWN"
$ws.Range("C40").Value = 40
$ws.Range("D40").Value = "11.02.18"
$ws.Range("E40").Value = "H.BROWN"
$ws.Range("F40").Value = 40
$ws.Range("G40").Value = "CACCTCC"
$ws.Range("H40").Value = "E7420L"

$ws.Range("A40").Interior.Color = $purple
$ws.Range("D40").Interior.Color = $purple
$ws.Range("G40").Interior.Color = $purple
$ws.Range("G40").WrapText = $true
$ws.Range("H40").WrapText = $true

$ws.Rows.Item(40).RowHeight = 16

# --- Selection / view state ----------------------------------------------
# Highlight the two newly-corrected index columns (B39:B40 and E39:E40).
$excel.Union($ws.Range("B39:B40"), $ws.Range("E39:E40")).Select()
$ws.Range("E39").Activate()
